# Correct misleading values in the "Dwellings_buildings" mapping scheme:
# revised area and cost assumptions for all occupancies; revised count
# assumptions for non-residential classifications.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2: H:1 / (was Mining and quarrying) -> All other industry
$ws.Range("B2").Value = "All other industry"
$ws.Range("C2").Value = 220
$ws.Range("D2").Value = 237.1016286644951

# Row 3: H:1 / (was Transport; storage and communication) -> Warehouses and storage
$ws.Range("B3").Value = "Warehouses and storage"
$ws.Range("C3").Value = 220
$ws.Range("D3").Value = 183.264332247557

# Row 4: H:1 / (was Manufacturing) -> Manufacturing and light industry
$ws.Range("B4").Value = "Manufacturing and light industry"
$ws.Range("C4").Value = 220
$ws.Range("D4").Value = 220.2434853420196

# Row 5: H:2 / (was Mining and quarrying) -> All other industry
$ws.Range("B5").Value = "All other industry"
$ws.Range("C5").Value = 480
$ws.Range("D5").Value = 237.1016286644951

# Row 6: H:2 / (was Transport; storage and communication) -> Warehouses and storage
$ws.Range("B6").Value = "Warehouses and storage"
$ws.Range("C6").Value = 480
$ws.Range("D6").Value = 183.264332247557

# Row 7: H:2 / (was Manufacturing) -> Manufacturing and light industry
$ws.Range("B7").Value = "Manufacturing and light industry"
$ws.Range("C7").Value = 480
$ws.Range("D7").Value = 220.2434853420196
